# Rename the first two season sheets to their zero-padded forms.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "2007-08"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "2008-09"

# Add the two missing game rows to the "2008-09" sheet (it previously only
# had the header row). Row/column layout mirrors every other season sheet:
#   A=GameId  B=Date  C=Home/Visitor flag  D=Team  E..H=quarter scores
#   I=Final score  J=money line
$ws2.Range("A2").Value = 1028
$ws2.Range("B2").Value = 39820
$ws2.Range("C2").Value = "V"
$ws2.Range("D2").Value = "Orlando"
$ws2.Range("E2").Value = 102
$ws2.Range("F2").Value = -42
$ws2.Range("G2").Value = 26
$ws2.Range("H2").Value = 20
$ws2.Range("I2").Value = 106
$ws2.Range("J2").Value = 120

$ws2.Range("A3").Value = 1029
$ws2.Range("B3").Value = 39820
$ws2.Range("C3").Value = "H"
$ws2.Range("D3").Value = "Atlanta"
$ws2.Range("E3").Value = 99
$ws2.Range("F3").Value = -54
$ws2.Range("G3").Value = 25
$ws2.Range("H3").Value = 32
$ws2.Range("I3").Value = 102
$ws2.Range("J3").Value = -140

# Match the formatting used by the analogous cells on every other season
# sheet: bold+bordered/centered style for the id column (A), and the
# custom date number format for the date column (B). Copy the formats
# from the first sheet, which already carries both styles, so we reuse
# the existing style records instead of minting new ones.
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)
$ws2.Range("A3").PasteSpecial(-4122)

$ws1.Range("B2").Copy()
$ws2.Range("B2").PasteSpecial(-4122)
$ws2.Range("B3").PasteSpecial(-4122)
